$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.101.76"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.814.39"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.46"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.02"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -4.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.809.13"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -8.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.713"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -8.20%  "
$ws.Range("E11").Value = "  -13.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000342"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -20.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.79"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -7.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.420.73"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.84"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -5.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.61"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +19.65%  "
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.808.82"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.34"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -6.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.271.09"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("E21").Value = "  -7.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.95"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -11.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.15"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -10.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.26"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -7.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.95"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.77"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +12.81%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "36.59"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -5.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.17"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -5.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -8.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "697.97"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.16"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.35"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.148"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -10.96%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.40"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -10.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.64"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0748"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0447"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -9.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("E43").Value = "  -10.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.41"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.67"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.09"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.01"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -7.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.03"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -4.95%  "
$ws.Range("E50").Value = "  -4.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -7.95%  "
